# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders / refreshes the worker account-statement rows (B16:G27):
#   - ADRIANA KARINA ALVEAR MARRUGO (1002203955) now occupies rows 16-22,
#     sorted by ascending "Periodo Mora" (1710..1804), with refreshed
#     Valor Mora / Salario Basico figures.
#   - JOHANA MEDINA MARRUGO (1051442443) now occupies rows 23-26, sorted
#     ascending (1806..1809), with a refreshed Salario Basico.
#   - JAIRO ANTONIO BELTRAN GUZMAN (92071331) moves from the first data
#     row down to the last one (row 27), values unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Doc = "CC"; NumDoc = "1002203955"; Nombre = "ADRIANA KARINA ALVEAR MARRUGO"; Periodo = "1710"; Mora = 32800;  Salario = 820000 },
    @{ Doc = "CC"; NumDoc = "1002203955"; Nombre = "ADRIANA KARINA ALVEAR MARRUGO"; Periodo = "1711"; Mora = 32800;  Salario = 820000 },
    @{ Doc = "CC"; NumDoc = "1002203955"; Nombre = "ADRIANA KARINA ALVEAR MARRUGO"; Periodo = "1712"; Mora = 32800;  Salario = 820000 },
    @{ Doc = "CC"; NumDoc = "1002203955"; Nombre = "ADRIANA KARINA ALVEAR MARRUGO"; Periodo = "1801"; Mora = 32800;  Salario = 820000 },
    @{ Doc = "CC"; NumDoc = "1002203955"; Nombre = "ADRIANA KARINA ALVEAR MARRUGO"; Periodo = "1802"; Mora = 32800;  Salario = 820000 },
    @{ Doc = "CC"; NumDoc = "1002203955"; Nombre = "ADRIANA KARINA ALVEAR MARRUGO"; Periodo = "1803"; Mora = 32800;  Salario = 820000 },
    @{ Doc = "CC"; NumDoc = "1002203955"; Nombre = "ADRIANA KARINA ALVEAR MARRUGO"; Periodo = "1804"; Mora = 1093;   Salario = 820000 },
    @{ Doc = "CC"; NumDoc = "1051442443"; Nombre = "JOHANA MEDINA MARRUGO";         Periodo = "1806"; Mora = 31249;  Salario = 1423500 },
    @{ Doc = "CC"; NumDoc = "1051442443"; Nombre = "JOHANA MEDINA MARRUGO";         Periodo = "1807"; Mora = 31249;  Salario = 1423500 },
    @{ Doc = "CC"; NumDoc = "1051442443"; Nombre = "JOHANA MEDINA MARRUGO";         Periodo = "1808"; Mora = 31249;  Salario = 1423500 },
    @{ Doc = "CC"; NumDoc = "1051442443"; Nombre = "JOHANA MEDINA MARRUGO";         Periodo = "1809"; Mora = 31249;  Salario = 1423500 },
    @{ Doc = "CC"; NumDoc = "92071331";   Nombre = "JAIRO ANTONIO BELTRAN GUZMAN";  Periodo = "2110"; Mora = 36341;  Salario = 908526 }
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.Doc
    $ws.Cells.Item($r, 3).Value = $row.NumDoc
    $ws.Cells.Item($r, 4).Value = $row.Nombre
    $ws.Cells.Item($r, 5).Value = $row.Periodo
    $ws.Cells.Item($r, 6).Value = $row.Mora
    $ws.Cells.Item($r, 7).Value = $row.Salario
    $r = $r + 1
}
